$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Update the constraint matrix cells
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1

# Activate the sheet and set the selection to match the saved view state
$ws.Activate()
$ws.Range("G5").Select()
